$d = $word.ActiveDocument
$d.Content.Find.Execute("generato da Jooq.", $true, $false, $false, $false, $false, $true, 1, $false, "generato da Jooq.", 2)
